$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 588.8
$ws.Range("I4").Value = 432.1111
$ws.Range("K4").Value = 432.1111
$ws.Range("M4").Value = -318.1111
$ws.Range("H11").Value = 324.6
$ws.Range("I11").Value = 324.6
$ws.Range("K11").Value = 324.6
$ws.Range("M11").Value = -184.6
$ws.Range("H17").Value = 1836.625
$ws.Range("J17").Value = 1815.5
$ws.Range("L17").Value = 5446.5
$ws.Range("N17").Value = -5782.5
$ws.Range("H76").Value = 1999
$ws.Range("J76").Value = 1999
$ws.Range("L76").Value = 1999
$ws.Range("N76").Value = -2629
$ws.Range("H79").Value = 1999
$ws.Range("J79").Value = 1999
$ws.Range("L79").Value = 1999
$ws.Range("N79").Value = -4183
$ws.Range("H106").Value = 3966.3333
$ws.Range("I106").Value = 3949.5
$ws.Range("K106").Value = 3949.5
$ws.Range("M106").Value = -3318.5
$ws.Range("H111").Value = 5428.7144
$ws.Range("I111").Value = 5833.5
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 17500.5
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -14433.5
$ws.Range("N111").Value = -15134
$ws.Range("H129").Value = 2203.7144
$ws.Range("I129").Value = 964.8889
$ws.Range("J129").Value = 4433.6
$ws.Range("K129").Value = 2894.6667
$ws.Range("L129").Value = 13300.8
$ws.Range("M129").Value = 2105.3333
$ws.Range("N129").Value = -23300.8
$ws.Range("H137").Value = 2181285.8
$ws.Range("I137").Value = 25000650
$ws.Range("J137").Value = 8012.905
$ws.Range("K137").Value = 75001950
$ws.Range("L137").Value = 24038.715
$ws.Range("M137").Value = -74999400
$ws.Range("N137").Value = -29138.715
$ws.Range("H138").Value = 4211.8237
$ws.Range("I138").Value = 8031
$ws.Range("J138").Value = 3702.6
$ws.Range("K138").Value = 24093
$ws.Range("L138").Value = 11107.8
$ws.Range("M138").Value = -18953
$ws.Range("N138").Value = -21387.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7060.9546
$ws.Range("I32").Value = 6407.8887
$ws.Range("J32").Value = 9999.75
$ws.Range("K32").Value = 6407.8887
$ws.Range("L32").Value = 9999.75
$ws.Range("M32").Value = -6120.8887
$ws.Range("N32").Value = -10573.75
$ws.Range("H61").Value = 3706.7917
$ws.Range("I61").Value = 2845.8333
$ws.Range("K61").Value = 2845.8333
$ws.Range("M61").Value = -2633.8333
$ws.Range("H102").Value = 5885.9165
$ws.Range("I102").Value = 5681.4443
$ws.Range("K102").Value = 5681.4443
$ws.Range("M102").Value = -4059.4443
$ws.Range("H132").Value = 5210838
$ws.Range("I132").Value = 2233.4707
$ws.Range("J132").Value = 25644596
$ws.Range("K132").Value = 6700.4121
$ws.Range("L132").Value = 76933788
$ws.Range("M132").Value = -4170.4121
$ws.Range("N132").Value = -76938848
$ws.Range("H136").Value = 3706.7917
$ws.Range("I136").Value = 2845.8333
$ws.Range("K136").Value = 8537.499899999999
$ws.Range("M136").Value = -5987.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17248198
$ws.Range("I20").Value = 25007980
$ws.Range("J20").Value = 4239
$ws.Range("K20").Value = 25007980
$ws.Range("L20").Value = 4239
$ws.Range("M20").Value = -25007733
$ws.Range("N20").Value = -4733
$ws.Range("H86").Value = 2201.6
$ws.Range("J86").Value = 3431
$ws.Range("L86").Value = 3431
$ws.Range("N86").Value = -5677
$ws.Range("H89").Value = 2201.6
$ws.Range("J89").Value = 3431
$ws.Range("L89").Value = 17155
$ws.Range("N89").Value = -28387
$ws.Range("H99").Value = 3094.3
$ws.Range("I99").Value = 2452.6316
$ws.Range("K99").Value = 2452.6316
$ws.Range("M99").Value = -954.6316000000002
$ws.Range("H107").Value = 1578.2858
$ws.Range("I107").Value = 1489.7646
$ws.Range("K107").Value = 1489.7646
$ws.Range("M107").Value = 430.2354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 91354.42999999999
$ws.Range("J52").Value = 91354.42999999999
$ws.Range("L52").Value = 91354.42999999999
$ws.Range("N52").Value = -91942.42999999999
$ws.Range("H58").Value = 2672.6
$ws.Range("I58").Value = 2201
$ws.Range("K58").Value = 2201
$ws.Range("M58").Value = -1998
$ws.Range("H60").Value = 12500.333
$ws.Range("I60").Value = 12199
$ws.Range("J60").Value = 12560.6
$ws.Range("K60").Value = 12199
$ws.Range("L60").Value = 12560.6
$ws.Range("M60").Value = -11688
$ws.Range("N60").Value = -13582.6
$ws.Range("H99").Value = 10349.2
$ws.Range("I99").Value = 11874.125
$ws.Range("K99").Value = 11874.125
$ws.Range("M99").Value = -10376.125
$ws.Range("H126").Value = 10349.2
$ws.Range("I126").Value = 11874.125
$ws.Range("K126").Value = 35622.375
$ws.Range("M126").Value = -33152.375
$ws.Range("H127").Value = 42695
$ws.Range("I127").Value = 30000
$ws.Range("J127").Value = 46926.668
$ws.Range("K127").Value = 30000
$ws.Range("L127").Value = 46926.668
$ws.Range("M127").Value = -25040
$ws.Range("N127").Value = -56846.668
$ws.Range("H136").Value = 2672.6
$ws.Range("I136").Value = 2201
$ws.Range("K136").Value = 6603
$ws.Range("M136").Value = -4053
$ws.Range("H139").Value = 64832.668
$ws.Range("J139").Value = 64832.668
$ws.Range("L139").Value = 64832.668
$ws.Range("N139").Value = -75112.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1471
$ws.Range("I46").Value = 1999
$ws.Range("J46").Value = 1295
$ws.Range("K46").Value = 5997
$ws.Range("L46").Value = 3885
$ws.Range("M46").Value = -5906
$ws.Range("N46").Value = -4067
$ws.Range("H58").Value = 1696
$ws.Range("I58").Value = 932.3333
$ws.Range("J58").Value = 3987
$ws.Range("K58").Value = 2796.9999
$ws.Range("L58").Value = 11961
$ws.Range("M58").Value = -2668.9999
$ws.Range("N58").Value = -12217
$ws.Range("H117").Value = 2031.5834
$ws.Range("J117").Value = 1908.7778
$ws.Range("L117").Value = 5726.3334
$ws.Range("N117").Value = -12610.3334
$ws.Range("H131").Value = 9763.223
$ws.Range("I131").Value = 17258.666
$ws.Range("K131").Value = 51775.99800000001
$ws.Range("M131").Value = -46735.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 150
$ws.Range("K2").Value = 150
$ws.Range("M2").Value = -37
$ws.Range("H70").Value = 73735.24000000001
$ws.Range("I70").Value = 128564.44
$ws.Range("J70").Value = 6253.154
$ws.Range("K70").Value = 128564.44
$ws.Range("L70").Value = 6253.154
$ws.Range("M70").Value = -128294.44
$ws.Range("N70").Value = -6793.154
$ws.Range("H73").Value = 73735.24000000001
$ws.Range("I73").Value = 128564.44
$ws.Range("J73").Value = 6253.154
$ws.Range("K73").Value = 128564.44
$ws.Range("L73").Value = 6253.154
$ws.Range("M73").Value = -127628.44
$ws.Range("N73").Value = -8125.154
$ws.Range("H132").Value = 2848.8518
$ws.Range("I132").Value = 1975.6111
$ws.Range("K132").Value = 5926.8333
$ws.Range("M132").Value = -3396.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 479.16
$ws.Range("I55").Value = 222
$ws.Range("K55").Value = 222
$ws.Range("M55").Value = -49
$ws.Range("H93").Value = 505.83334
$ws.Range("I93").Value = 505.83334
$ws.Range("K93").Value = 505.83334
$ws.Range("M93").Value = 742.16666
$ws.Range("H100").Value = 1532.8889
$ws.Range("I100").Value = 1513.8572
$ws.Range("K100").Value = 1513.8572
$ws.Range("M100").Value = -972.8571999999999
$ws.Range("H120").Value = 199998.17
$ws.Range("J120").Value = 199998.17
$ws.Range("L120").Value = 199998.17
$ws.Range("N120").Value = -209674.17
$ws.Range("H127").Value = 49999.5
$ws.Range("J127").Value = 49999.5
$ws.Range("L127").Value = 49999.5
$ws.Range("N127").Value = -59919.5
$ws.Range("H132").Value = 4436.143
$ws.Range("I132").Value = 3041.4285
$ws.Range("J132").Value = 6528.2144
$ws.Range("K132").Value = 9124.2855
$ws.Range("L132").Value = 19584.6432
$ws.Range("M132").Value = -6594.2855
$ws.Range("N132").Value = -24644.6432
$ws.Range("H137").Value = 53315.105
$ws.Range("J137").Value = 55554.832
$ws.Range("L137").Value = 55554.832
$ws.Range("N137").Value = -65754.83199999999
$ws.Range("H139").Value = 63701.332
$ws.Range("I139").Value = 50650
$ws.Range("J139").Value = 70227
$ws.Range("K139").Value = 50650
$ws.Range("L139").Value = 70227
$ws.Range("M139").Value = -45510
$ws.Range("N139").Value = -80507

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 22524
$ws.Range("I55").Value = 22524
$ws.Range("K55").Value = 22524
$ws.Range("M55").Value = -22247
$ws.Range("H56").Value = 41203.125
$ws.Range("J56").Value = 46375
$ws.Range("L56").Value = 46375
$ws.Range("N56").Value = -47803
$ws.Range("H122").Value = 20836300
$ws.Range("I122").Value = 2677.2222
$ws.Range("K122").Value = 8031.6666
$ws.Range("M122").Value = -5581.6666
$ws.Range("H124").Value = 52500
$ws.Range("J124").Value = 52500
$ws.Range("L124").Value = 52500
$ws.Range("N124").Value = -62320
$ws.Range("H133").Value = 88377.39999999999
$ws.Range("J133").Value = 88377.39999999999
$ws.Range("L133").Value = 88377.39999999999
$ws.Range("N133").Value = -98497.39999999999
